$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63, shifting existing rows 63-70 down to 64-71
$ws.Rows.Item(63).Insert()

$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44769
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100112029
$ws.Cells.Item(63, 7).Value = "Orégano"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 16
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 13).Value = 20000
$ws.Cells.Item(63, 14).Value = "$/docena de atados"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 6667
$ws.Cells.Item(63, 17).Value = 3
$ws.Cells.Item(63, 18).Value = "Hortaliza"
